$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '28.278.22'
$ws.Range("E2").Value = '  -2.68%  '
$ws.Range("D3").Value = '1.935.40'
$ws.Range("E3").Value = '  -1.38%  '
$ws.Range("D4").Value = '1.013'
$ws.Range("E4").Value = '  +0.80%  '
$ws.Range("D5").Value = '321.16'
$ws.Range("E5").Value = '  -1.87%  '
$ws.Range("D6").Value = '1.011'
$ws.Range("E6").Value = '  +0.67%  '
$ws.Range("D7").Value = '0.4752'
$ws.Range("E7").Value = '  -4.50%  '
$ws.Range("D8").Value = '0.4060'
$ws.Range("E8").Value = '  -3.72%  '
$ws.Range("D9").Value = '53.52'
$ws.Range("E9").Value = '  +1.56%  '
$ws.Range("D10").Value = '0.08502'
$ws.Range("E10").Value = '  -7.75%  '
$ws.Range("D11").Value = '1.052'
$ws.Range("E11").Value = '  -4.21%  '
$ws.Range("D12").Value = '22.31'
$ws.Range("E12").Value = '  -2.42%  '
$ws.Range("D13").Value = '1.937.34'
$ws.Range("E13").Value = '  -2.37%  '
$ws.Range("D14").Value = '7.528'
$ws.Range("E14").Value = '  -4.00%  '
$ws.Range("D15").Value = '6.128'
$ws.Range("E15").Value = '  -4.94%  '
$ws.Range("D16").Value = '1.014'
$ws.Range("E16").Value = '  +0.75%  '
$ws.Range("D17").Value = '89.97'
$ws.Range("E17").Value = '  -1.78%  '
$ws.Range("D18").Value = '0.00001071'
$ws.Range("E18").Value = '  -2.54%  '
$ws.Range("D19").Value = '0.06615'
$ws.Range("E19").Value = '  -1.02%  '
$ws.Range("D20").Value = '18.23'
$ws.Range("E20").Value = '  -5.42%  '
$ws.Range("D21").Value = '1.011'
$ws.Range("E21").Value = '  +0.68%  '
$ws.Range("D22").Value = '5.801'
$ws.Range("E22").Value = '  -2.54%  '
$ws.Range("D23").Value = '28.304.83'
$ws.Range("E23").Value = '  -2.67%  '
$ws.Range("D24").Value = '11.44'
$ws.Range("E24").Value = '  -4.93%  '
$ws.Range("D25").Value = '2.317'
$ws.Range("E25").Value = '  +1.75%  '
$ws.Range("D26").Value = '2.165.41'
$ws.Range("E26").Value = '  -2.26%  '
$ws.Range("D27").Value = '155.31'
$ws.Range("E27").Value = '  -0.54%  '
$ws.Range("D28").Value = '20.23'
$ws.Range("E28").Value = '  -1.82%  '
$ws.Range("D29").Value = '2.173'
$ws.Range("E29").Value = '  -3.67%  '
$ws.Range("D30").Value = '5.783'
$ws.Range("E30").Value = '  -8.19%  '
$ws.Range("D31").Value = '124.04'
$ws.Range("E31").Value = '  -1.73%  '
$ws.Range("D32").Value = '0.9848'
$ws.Range("E32").Value = '  -5.83%  '
$ws.Range("D33").Value = '0.09621'
$ws.Range("E33").Value = '  -2.13%  '
$ws.Range("D34").Value = '1.445'
$ws.Range("E34").Value = '  -5.45%  '
$ws.Range("D35").Value = '3.667'
$ws.Range("E35").Value = '  -0.41%  '
$ws.Range("D36").Value = '5.597'
$ws.Range("E36").Value = '  -3.70%  '
$ws.Range("D37").Value = '9.280'
$ws.Range("E37").Value = '  +2.94%  '
$ws.Range("D38").Value = '0.02323'
$ws.Range("E38").Value = '  -4.15%  '
$ws.Range("D39").Value = '0.06187'
$ws.Range("E39").Value = '  -2.76%  '
$ws.Range("D40").Value = '1.243'
$ws.Range("E40").Value = '  -5.59%  '
$ws.Range("D41").Value = '0.6205'
$ws.Range("E41").Value = '  -3.63%  '
$ws.Range("D42").Value = '11.16'
$ws.Range("E42").Value = '  -2.23%  '
$ws.Range("D43").Value = '1.011'
$ws.Range("E43").Value = '  +0.62%  '
$ws.Range("D44").Value = '0.1913'
$ws.Range("E44").Value = '  -3.38%  '
$ws.Range("D45").Value = '1.327'
$ws.Range("E45").Value = '  +0.39%  '
$ws.Range("D46").Value = '0.5926'
$ws.Range("E46").Value = '  -4.91%  '
$ws.Range("D47").Value = '12.86'
$ws.Range("E47").Value = '  -3.83%  '
$ws.Range("D48").Value = '2.050'
$ws.Range("E48").Value = '  -6.90%  '
$ws.Range("D49").Value = '3.397'
$ws.Range("E49").Value = '  -2.03%  '
$ws.Range("D50").Value = '0.06789'
$ws.Range("E50").Value = '  -2.99%  '
$ws.Range("D51").Value = '110.14'
$ws.Range("E51").Value = '  -1.81%  '
